# Y4_B2526_Excuses.xlsx - attendance app re-upload
# Updates the two excuse log rows: new student IDs, new log date, and
# row 3's log time switches from a numeric time value to a plain text
# time string (matching row 2's existing "Log Time" formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 -----------------------------------------------------------
# A2: Student ID 211276 -> 201252 (kept as text, same cell style)
$ws.Range("A2").Formula = '="201252"'
$ws.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4163)   # xlPasteValues: flatten formula -> literal text

# C2: Log Date 16/10/2025 -> 18/10/2025 (already plain text)
$ws.Range("C2").Value = "18/10/2025"

# --- Row 3 -----------------------------------------------------------
# A3: Student ID 212155 -> 201572 (kept as text, same cell style)
$ws.Range("A3").Formula = '="201572"'
$ws.Range("A3").Copy()
$ws.Range("A3").PasteSpecial(-4163)   # xlPasteValues: flatten formula -> literal text

# C3: Log Date 16/10/2025 -> 18/10/2025 (already plain text)
$ws.Range("C3").Value = "18/10/2025"

# D3: Log Time changes from numeric time-of-day (0.4375) to plain text
# "10:30:00", matching the formatting already used by D2/E3/etc.
$ws.Range("D3").Formula = '="10:30:00"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)   # xlPasteValues: flatten formula -> literal text
$ws.Range("E3").Copy()
$ws.Range("D3").PasteSpecial(-4122)   # xlPasteFormats: adopt the row's plain-text style
